$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 previously held the category "sloveso" (verb); rename it to "verbs"
$ws.Cells.Item(4, 3).Value2 = "verbs"

# Rows 9-50 previously had no category in column C at all; give them all
# the new "verbs" category (mirrors rows 9-26 which already said "sloveso"
# plus rows 27-50 which had nothing).
for ($r = 9; $r -le 50; $r++) {
    $ws.Cells.Item($r, 3).Value2 = "verbs"
}

# Move the active selection down to B51, matching the author's new cursor
# position after adding the new rows/column of data.
$ws.Range("B51").Select()
